$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set D5:D7 to "Yes" (these cells have list validation "Yes,No")
$ws.Range("D5").Value = "Yes"
$ws.Range("D6").Value = "Yes"
$ws.Range("D7").Value = "Yes"

# Update the selection to D7 as recorded in the sheet view
$ws.Range("D7").Select()
